# Regenerate save_data to use K instead of Strike# in column G.
# Writes the newly computed K values (std/mean recalculated, s_vals written)
# into column G ("K" header) for each row of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 3
    5  = 1
    6  = 1
    7  = 4
    8  = 3
    9  = 4
    10 = 5
    11 = 6
    12 = 2
    13 = 7
    14 = 3
    15 = 8
    16 = 4
    17 = 5
    18 = 6
    19 = 3
    20 = 5
    21 = 2
    22 = 3
    24 = 0
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
